$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

function Set-PlainValue($ref, $val) {
    $ws.Range($ref).Value = $val
}

Set-TextValue "D2" "27.597.42"
Set-PlainValue "E2" "  -1.75%  "
Set-TextValue "D3" "1.666.69"
Set-PlainValue "E3" "  -3.51%  "
Set-PlainValue "E4" "  +0.06%  "
Set-TextValue "D5" "215.46"
Set-PlainValue "E5" "  -1.52%  "
Set-PlainValue "E6" "  -2.17%  "
Set-PlainValue "E7" "  +0.03%  "
Set-TextValue "D8" "23.65"
Set-PlainValue "E8" "  -1.76%  "
Set-TextValue "D9" "0.263"
Set-PlainValue "E9" "  -0.42%  "
Set-PlainValue "E10" "  -2.02%  "
Set-TextValue "D11" "0.0883"
Set-TextValue "D12" "1.901.61"
Set-PlainValue "E12" "  -3.56%  "
Set-TextValue "D13" "1.664.53"
Set-PlainValue "E13" "  -3.60%  "
Set-PlainValue "E14" "  -2.58%  "
Set-TextValue "D15" "0.559"
Set-PlainValue "E15" "  -0.76%  "
Set-TextValue "D16" "66.18"
Set-PlainValue "E16" "  -2.29%  "
Set-TextValue "D17" "245.93"
Set-PlainValue "E17" "  +0.91%  "
Set-TextValue "D18" "27.620.85"
Set-PlainValue "E18" "  -1.43%  "
Set-PlainValue "E19" "  -3.56%  "
Set-TextValue "D20" "7.55"
Set-PlainValue "E20" "  -4.22%  "
Set-PlainValue "E21" "  +0.00%  "
Set-PlainValue "E22" "  -3.23%  "
Set-TextValue "D23" "9.30"
Set-PlainValue "E23" "  -4.82%  "
Set-TextValue "D24" "2.05"
Set-PlainValue "E24" "  -3.98%  "
Set-TextValue "D25" "146.28"
Set-PlainValue "E25" "  -1.95%  "
Set-TextValue "D26" "7.19"
Set-PlainValue "E26" "  -4.54%  "
Set-TextValue "D27" "16.40"
Set-PlainValue "E27" "  -2.35%  "
Set-PlainValue "E28" "  +0.11%  "
Set-PlainValue "E29" "  -2.81%  "
Set-PlainValue "E30" "  +3.58%  "
Set-TextValue "D31" "0.0505"
Set-PlainValue "E31" "  -1.09%  "
Set-PlainValue "E32" "  -2.83%  "
Set-TextValue "D33" "1.478.63"
Set-PlainValue "E33" "  -0.77%  "
Set-TextValue "D34" "3.12"
Set-PlainValue "E34" "  -5.09%  "
Set-PlainValue "E35" "  -5.86%  "
Set-TextValue "D36" "0.938"
Set-PlainValue "E36" "  -2.31%  "
Set-PlainValue "E37" "  -1.13%  "
Set-PlainValue "B38" "VeChain"
Set-PlainValue "C38" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D38" "0.0171"
Set-PlainValue "E38" "  -2.54%  "
Set-PlainValue "B39" "ImmutableX"
Set-PlainValue "C39" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D39" "0.572"
Set-PlainValue "E39" "  -6.34%  "
Set-TextValue "D40" "69.69"
Set-PlainValue "E40" "  -2.07%  "
Set-PlainValue "E41" "  -5.80%  "
Set-TextValue "D42" "0.999"
Set-PlainValue "E42" "  +0.02%  "
Set-TextValue "D43" "5.42"
Set-PlainValue "E43" "  -7.23%  "
Set-TextValue "D44" "1.809.80"
Set-PlainValue "E44" "  -3.50%  "
Set-PlainValue "E45" "  -4.16%  "
Set-PlainValue "E46" "  -1.07%  "
Set-PlainValue "E47" "  -2.77%  "
Set-TextValue "D48" "89.33"
Set-PlainValue "E48" "  -2.57%  "
Set-PlainValue "E49" "  -1.83%  "
Set-PlainValue "E50" "  -2.89%  "
Set-TextValue "D51" "7.94"
Set-PlainValue "E51" "  -3.57%  "
